$wb = $excel.ActiveWorkbook

# --- Tree_Sum sheet: add follow-up notes next to Path Sum / Sum Root to Leaf Numbers ---
$wsSum = $wb.Worksheets.Item("Tree_Sum")
$wsSum.Range("B3").Value = "//为什么要先把root加进去？为什么上面的只一行的方法不行？"
$wsSum.Range("B4").Value = "算法的本质是一次先序遍历（为啥？）"
$wsSum.Range("B4").Select()

# --- Tree_Traversal sheet: add a note, becomes the active sheet ---
$wsTrav = $wb.Worksheets.Item("Tree_Traversal")
$wsTrav.Range("B2").Value = "//item.clear();【注】，这样写错误！原因问问老师。"
$wsTrav.Range("B2").NumberFormat = "d-mmm-yy"
$wsTrav.Activate()
$wsTrav.Range("B2").Select()
